$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right before the current row 492, shifting old rows 492:570 down to 496:574.
$ws.Rows("492:495").Insert()

# New row 492 data
$ws.Cells.Item(492, 1).Value = 7
$ws.Cells.Item(492, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(492, 3).Value = "Ñuble"
$ws.Cells.Item(492, 4).Value = "2023-09-11"
$ws.Cells.Item(492, 5).Value = 16
$ws.Cells.Item(492, 6).Value = 100112008
$ws.Cells.Item(492, 7).Value = "Coliflor"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 300
$ws.Cells.Item(492, 11).Value = 900
$ws.Cells.Item(492, 12).Value = 900
$ws.Cells.Item(492, 13).Value = 900
$ws.Cells.Item(492, 14).Value = "`$/unidad"
$ws.Cells.Item(492, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(492, 16).Value = 900
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"

# New row 493 data
$ws.Cells.Item(493, 1).Value = 7
$ws.Cells.Item(493, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(493, 3).Value = "Ñuble"
$ws.Cells.Item(493, 4).Value = "2023-09-11"
$ws.Cells.Item(493, 5).Value = 16
$ws.Cells.Item(493, 6).Value = 100112008
$ws.Cells.Item(493, 7).Value = "Coliflor"
$ws.Cells.Item(493, 8).Value = "Sin especificar"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 300
$ws.Cells.Item(493, 11).Value = 1000
$ws.Cells.Item(493, 12).Value = 1000
$ws.Cells.Item(493, 13).Value = 1000
$ws.Cells.Item(493, 14).Value = "`$/unidad"
$ws.Cells.Item(493, 15).Value = "Región del Maule"
$ws.Cells.Item(493, 16).Value = 1000
$ws.Cells.Item(493, 17).Value = 1
$ws.Cells.Item(493, 18).Value = "Hortaliza"

# New row 494 data
$ws.Cells.Item(494, 1).Value = 7
$ws.Cells.Item(494, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(494, 3).Value = "Ñuble"
$ws.Cells.Item(494, 4).Value = "2023-09-11"
$ws.Cells.Item(494, 5).Value = 16
$ws.Cells.Item(494, 6).Value = 100112008
$ws.Cells.Item(494, 7).Value = "Coliflor"
$ws.Cells.Item(494, 8).Value = "Sin especificar"
$ws.Cells.Item(494, 9).Value = "Segunda"
$ws.Cells.Item(494, 10).Value = 300
$ws.Cells.Item(494, 11).Value = 700
$ws.Cells.Item(494, 12).Value = 700
$ws.Cells.Item(494, 13).Value = 700
$ws.Cells.Item(494, 14).Value = "`$/unidad"
$ws.Cells.Item(494, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(494, 16).Value = 700
$ws.Cells.Item(494, 17).Value = 1
$ws.Cells.Item(494, 18).Value = "Hortaliza"

# New row 495 data
$ws.Cells.Item(495, 1).Value = 7
$ws.Cells.Item(495, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(495, 3).Value = "Ñuble"
$ws.Cells.Item(495, 4).Value = "2023-09-11"
$ws.Cells.Item(495, 5).Value = 16
$ws.Cells.Item(495, 6).Value = 100112008
$ws.Cells.Item(495, 7).Value = "Coliflor"
$ws.Cells.Item(495, 8).Value = "Sin especificar"
$ws.Cells.Item(495, 9).Value = "Segunda"
$ws.Cells.Item(495, 10).Value = 300
$ws.Cells.Item(495, 11).Value = 800
$ws.Cells.Item(495, 12).Value = 800
$ws.Cells.Item(495, 13).Value = 800
$ws.Cells.Item(495, 14).Value = "`$/unidad"
$ws.Cells.Item(495, 15).Value = "Región del Maule"
$ws.Cells.Item(495, 16).Value = 800
$ws.Cells.Item(495, 17).Value = 1
$ws.Cells.Item(495, 18).Value = "Hortaliza"
